# "new format for excel"
#
# The guild-config sheet used to store its six "on/off" columns (B:G) for
# rows 3-9 as TRUE/FALSE booleans. The new format stores the very same
# 0 values, but as plain numbers instead of booleans (no more t="b" cells
# in the exported XML) - existing cell styles/borders/fills are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 3-9, columns B-G: rewrite the boolean FALSE cells as numeric 0.
$ws.Range("B3:G9").Value = 0
